$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 46, shifting existing rows 46-92 down to 47-93
$ws.Rows("46:46").Insert()

# Populate the newly inserted row 46 with the new observation
$ws.Cells.Item(46, 1).Value = 1
$ws.Cells.Item(46, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(46, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(46, 4).Value = 44778
$ws.Cells.Item(46, 5).Value = 15
$ws.Cells.Item(46, 6).Value = 100112021
$ws.Cells.Item(46, 7).Value = "Ají"
$ws.Cells.Item(46, 8).Value = "Inferno"
$ws.Cells.Item(46, 9).Value = "Primera"
$ws.Cells.Item(46, 10).Value = 160
$ws.Cells.Item(46, 11).Value = 11000
$ws.Cells.Item(46, 12).Value = 12000
$ws.Cells.Item(46, 13).Value = 11500
$ws.Cells.Item(46, 14).Value = "$/caja 15 kilos"
$ws.Cells.Item(46, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(46, 16).Value = 767
$ws.Cells.Item(46, 17).Value = 15
$ws.Cells.Item(46, 18).Value = "Hortaliza"
